# Add "Variadic templates" slide (new slide 4) before the existing
# "What was left out?" slide (which becomes slide 5), and tweak a few
# shape positions on the "Function templates" slide (slide 3).

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 3 ("Function templates"): nudge a few shapes.
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

$tb5 = $s3.Shapes.Item(3)   # "TextBox 5"
$tb5.Left   = 37.69347456692914
$tb5.Top    = 132.50544307086614
$tb5.Width  = 322.3065454330708
$tb5.Height = 104.20780527559056

$tb7 = $s3.Shapes.Item(4)   # "TextBox 7"
$tb7.Left   = 252.43552181102362
$tb7.Top    = 179.05670291338583
$tb7.Width  = 377.80638795275587
$tb7.Height = 104.20780527559056

$grp15 = $s3.Shapes.Item(5) # "Group 15"
$grp15.Left = 582.6017422834645
$grp15.Top  = 231.69945881889765

$grp14 = $s3.Shapes.Item(6) # "Group 14"
$grp14.Left = 95.91536433070867
$grp14.Top  = 334.445915511811

# ---------------------------------------------------------------------
# 2) Insert the new "Variadic templates" slide before "What was left
#    out?" by duplicating that slide and moving the copy one slot
#    earlier; then re-purpose the moved copy's content.
# ---------------------------------------------------------------------
$oldSlide = $p.Slides.Item(4)          # "What was left out?"
$dupRange = $oldSlide.Duplicate()
$newSlide = $dupRange.Item(1)
$newSlide.MoveTo(4)

# --- 2a) Title ---
$title = $newSlide.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Variadic"
$title.TextFrame.TextRange.InsertAfter(" templates") | Out-Null

# --- 2b) Content placeholder ---
$content = $newSlide.Shapes.Item(2)
$content.Name = "Content Placeholder 4"
$content.TextFrame.TextRange.Text = "Implementing function with arbitrary number of arguments"

# --- 2c) New code textbox ---
$code = $newSlide.Shapes.AddTextbox(1, 81.04835645669291, 223.08615173228347, 534.4356005511811, 162.37032496062992)
$code.Name = "TextBox 3"
$code.Fill.ForeColor.RGB = 3156734
$code.Line.ForeColor.RGB = 0

$codeTr = $code.TextFrame.TextRange
$codeTr.Text = "double sum() { return 0.0; }`rtemplate<typename T, typename... Tail>`rdouble sum(T head, Tail... tail) {`r    return head + sum(tail...);`r}`r…`rstd::cout << sum(1.2, 2.3, 3.4) << std::endl;`rstd::cout << sum(1.2, 2.3, 3.4, 4.5) << std::endl;"
$codeTr.Font.Name = "Courier New"
$codeTr.Font.Size = 16

# ---------------------------------------------------------------------
# 3) Old "What was left out?" slide now sits at position 5 - add the
#    missing bullet about container templates. The slide-number field
#    auto-updates to "5" because of its new position.
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$s5content = $s5.Shapes.Item(2)
$s5content.TextFrame.TextRange.Text = "Container templates, i.e., writing your own generic containers"
